$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.360.83'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '3.142.02'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = "'609.19"
$ws.Range("D6").Value = "'143.80"
$ws.Range("E6").Value = '  -2.54%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '3.141.99'
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D9").Value = "'0.531"
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("D10").Value = "'0.151"
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("D11").Value = "'5.38"
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("D12").Value = "'0.473"
$ws.Range("E12").Value = '  -0.80%  '
$ws.Range("E13").Value = '  +1.84%  '
$ws.Range("D14").Value = "'35.45"
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("D15").Value = '3.657.75'
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("E16").Value = '  +2.51%  '
$ws.Range("D17").Value = '64.326.05'
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").Value = '3.142.67'
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").Value = "'6.88"
$ws.Range("E19").Value = '  -0.88%  '
$ws.Range("D20").Value = "'477.24"
$ws.Range("E20").Value = '  -0.79%  '
$ws.Range("D21").Value = "'14.80"
$ws.Range("E21").Value = '  +0.74%  '
$ws.Range("D22").Value = "'0.720"
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("D23").Value = "'7.81"
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("D24").Value = "'85.71"
$ws.Range("E24").Value = '  +2.48%  '
$ws.Range("E25").Value = '  -0.41%  '
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("D27").Value = "'2.77"
$ws.Range("E27").Value = '  -3.53%  '
$ws.Range("D28").Value = "'8.46"
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").Value = "'7.40"
$ws.Range("E29").Value = '  +9.99%  '
$ws.Range("D30").Value = "'0.116"
$ws.Range("E30").Value = '  +3.31%  '
$ws.Range("E31").Value = '  -5.79%  '
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").Value = "'26.70"
$ws.Range("E33").Value = '  +1.77%  '
$ws.Range("E34").Value = '  -3.61%  '
$ws.Range("E35").Value = '  +0.71%  '
$ws.Range("D36").Value = "'5.97"
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").Value = "'52.77"
$ws.Range("E37").Value = '  -3.24%  '
$ws.Range("D38").Value = '0.0₃0744'
$ws.Range("E38").Value = '  +2.69%  '
$ws.Range("D39").Value = "'450.13"
$ws.Range("E39").Value = '  -0.58%  '
$ws.Range("D40").Value = "'2.98"
$ws.Range("E40").Value = '  +1.49%  '
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").Value = "'0.118"
$ws.Range("E42").Value = '  +0.57%  '
$ws.Range("E43").Value = '  -1.16%  '
$ws.Range("D44").Value = '2.880.11'
$ws.Range("E44").Value = '  +1.32%  '
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("E46").Value = '  -1.39%  '
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").Value = "'34.24"
$ws.Range("E51").Value = '  +7.22%  '
